# Generate Report for Handback
# Replaces the file-handback identifiers (UUID-based file names / hash / timestamps)
# recorded on the "Overview", "zh-cn" and "de-de" sheets with the values produced by
# the latest handback run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rows 2 and 3 reference the two handed-back files by name.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsOverview.Range("B2").Value = "e2e\6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsOverview.Range("G2").Value = "2016-08-26 21:01:13"

$wsOverview.Range("A3").Value = "ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsOverview.Range("B3").Value = "e2e\ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsOverview.Range("G3").Value = "2016-08-26 21:01:13"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": source file names, generated xliff name and the handoff /
# handback datetimes.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsZhCn.Range("G2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 21:01:02"
$wsZhCn.Range("I2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsZhCn.Range("J2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-26 21:01:31"

$wsZhCn.Range("A3").Value = "ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsZhCn.Range("G3").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 21:01:02"
$wsZhCn.Range("I3").Value = "ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsZhCn.Range("J3").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-26 21:01:31"

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape as "zh-cn" but for the de-de xliff.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsDeDe.Range("G2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 21:01:13"
$wsDeDe.Range("I2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.md"
$wsDeDe.Range("J2").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-26 21:01:38"

$wsDeDe.Range("A3").Value = "ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsDeDe.Range("G3").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 21:01:13"
$wsDeDe.Range("I3").Value = "ffff38868c6b-a99a-4dc2-a3f7-8291c9cd6965.md"
$wsDeDe.Range("J3").Value = "6b42df40-8730-4109-a8cc-654bc1314007.d5178edfda8971a529879990027df4de566ee28b.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-26 21:01:38"
